# Update 1.3: Added Object Spenditure Report and other
#
# The author inserted a new "Бригадир" (foreman) column between the
# existing "Бригада" and "Дата" columns (i.e. before the old column E),
# shifting the remaining header columns one place to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the existing header values from E1:J1 before shifting them one
# column to the right, to make room for the new "Бригадир" column. Reading
# into variables first avoids clobbering values we still need to copy.
$vE = $ws.Range("E1").Value()
$vF = $ws.Range("F1").Value()
$vG = $ws.Range("G1").Value()
$vH = $ws.Range("H1").Value()
$vI = $ws.Range("I1").Value()
$vJ = $ws.Range("J1").Value()

$ws.Range("K1").Value = $vJ
$ws.Range("K1").Font().Bold = $true
$ws.Range("J1").Value = $vI
$ws.Range("I1").Value = $vH
$ws.Range("H1").Value = $vG
$ws.Range("G1").Value = $vF
$ws.Range("F1").Value = $vE
$ws.Range("E1").Value = "Бригадир"
$ws.Range("E1").Font().Bold = $true

# Re-fit the columns whose header text changed or moved so that the
# column widths stay readable for the new layout.
$ws.Columns("C:C").ColumnWidth = 10.833333333333332
$ws.Columns("E:E").ColumnWidth = 8.666666666666666
$ws.Columns("F:F").ColumnWidth = 14.0
$ws.Columns("G:G").ColumnWidth = 13.666666666666666
$ws.Columns("J:J").ColumnWidth = 11.666666666666666

# Restore the active selection left by the author on the sheet.
[void]$ws.Range("J7").Select()
